$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they keep their exact
# textual representation (e.g. trailing zeros) instead of becoming a Number
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '35.327.89'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '1.899.23'
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '243.73'
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("D6").Value = '0.650'
$ws.Range("E6").Value = '  +4.89%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '41.52'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = '0.341'
$ws.Range("E9").Value = '  +4.52%  '
$ws.Range("D10").Value = '50.12'
$ws.Range("E10").Value = '  +7.82%  '
$ws.Range("D11").Value = '0.0710'
$ws.Range("E11").Value = '  +2.93%  '
$ws.Range("D12").Value = '0.0996'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("D14").Value = '12.11'
$ws.Range("E14").Value = '  +6.63%  '
$ws.Range("E15").Value = '  +3.01%  '
$ws.Range("D16").Value = '1.900.23'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = '4.84'
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("D18").Value = '35.375.70'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").Value = '71.59'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").Value = '0.0₃0815'
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").Value = '241.89'
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").Value = '12.50'
$ws.Range("E22").Value = '  +3.14%  '
$ws.Range("D23").Value = '4.73'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D26").Value = '2.35'
$ws.Range("E26").Value = '  +28.14%  '
$ws.Range("D27").Value = '170.16'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '8.33'
$ws.Range("E28").Value = '  +4.72%  '
$ws.Range("D29").Value = '18.19'
$ws.Range("E29").Value = '  +3.66%  '
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("D31").Value = '4.12'
$ws.Range("E31").Value = '  +3.62%  '
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").Value = '1.01'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '0.927'
$ws.Range("E34").Value = '  +19.31%  '
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("E36").Value = '  +3.43%  '
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("D38").Value = '1.33'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").Value = '0.0210'
$ws.Range("E39").Value = '  +4.25%  '
$ws.Range("E40").Value = '  +1.99%  '
$ws.Range("D41").Value = '0.0629'
$ws.Range("E41").Value = '  +14.07%  '
$ws.Range("D42").Value = '15.80'
$ws.Range("E42").Value = '  +5.90%  '
$ws.Range("D43").Value = '89.51'
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").Value = '1.338.11'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").Value = '47.46'
$ws.Range("E46").Value = '  +38.72%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("E48").Value = '  +1.87%  '
$ws.Range("D49").Value = '12.38'
$ws.Range("E49").Value = '  -8.56%  '
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '2.082.25'
$ws.Range("E51").Value = '  +2.49%  '
